$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A131").Value = "grape juice"
$ws.Range("A132").Value = "orange juice"
$ws.Range("A133").Value = "grape seeds"
